# results.xlsx refresh: "change nested loop to if gen and replace
# if(xy[i]) o +=xy[i] in hamming.v"
#
# The hamming-weight generator change alters the Synopsys XOR/IV counts
# (columns P, Q, S) for the 8-bit ("Size"=128/256 family -> row 5, here the
# 2nd data block) and 16-bit rows of the results sheet; R/T/U are formulas
# and recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (first results row): new XOR / IV / NOR-weight counts ---
$ws.Range("P5").Value = 496
$ws.Range("Q5").Value = 332
$ws.Range("S5").Value = 486

# --- Row 6 (second results row): new XOR / IV / NOR-weight counts ---
$ws.Range("P6").Value = 4825
$ws.Range("Q6").Value = 3126
$ws.Range("S6").Value = 4810

# --- Recalculate so the SUM()-based formula cells (R5,T5,U5,R6,T6,U6) pick
#     up the new cached values right away. ---
$excel.Calculate()

# --- View-state bookkeeping matching the author's last on-screen position ---
$wb.Windows.Item(1).TabRatio = 151
$ws.Range("L1").Select() | Out-Null
$ws.Range("V6").Select() | Out-Null

# --- Header style (P1:U1) should carry the "apply border" flag alongside
#     its existing font/alignment overrides. ---
$ws.Range("P1:U1").Style.IncludeBorder = $true
